$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 68: appointments - University Assistant Professor (Beginning Sep. 2022) - Institute of Criminology, University of Cambridge - 2022
$ws.Range("A68").Value = "appointments"
$ws.Range("B68").Value = "University Assistant Professor (Beginning Sep. 2022) "
$ws.Range("C68").Value = "2022"
$ws.Range("E68").Value = "Institute of Criminology, University of Cambridge"

# Row 69: work - Associate - Department of Sociology, Harvard University - Apr 2022 - Present
$ws.Range("A69").Value = "work"
$ws.Range("B69").Value = "Associate"
$ws.Range("C69").Value = "Apr 2022"
$ws.Range("D69").Value = "Present"
$ws.Range("E69").Value = "Department of Sociology, Harvard University"
